# Minor fixes on DeveloperGuide and comments (#240)
# Update the Logic component class diagram: "AddCommand" -> "AddmedsCommand",
# "FindCommand" -> "ViewCommand", and shrink the example text in the
# "comment" note (Folded Corner shape) from 10.5pt to 10pt, matching the
# re-wording of Address Book -> HealthBase commands elsewhere in the guide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the folded-corner "comment" shape that documents example command
# names (XYZCommand = AddCommand, FindCommand, etc.)
$shape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.Name -eq "Folded Corner 126") {
        $shape = $candidate
    }
}

# Replace "AddCommand" with "AddmedsCommand" (found dynamically so the
# fixed-up text of the previous replacement doesn't shift later offsets).
$text = $shape.TextFrame.TextRange.Text
$pos = $text.IndexOf("AddCommand")
$range = $shape.TextFrame.TextRange.Characters($pos + 1, 10)
$range.Text = "AddmedsCommand"

# Replace "FindCommand" with "ViewCommand".
$text = $shape.TextFrame.TextRange.Text
$pos = $text.IndexOf("FindCommand")
$range = $shape.TextFrame.TextRange.Characters($pos + 1, 11)
$range.Text = "ViewCommand"

# The whole note shrinks from 10.5pt to 10pt as part of the edit.
$shape.TextFrame.TextRange.Font.Size = 10
